$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.034.12'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.706.17'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.39'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4000'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4044'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.475'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.37'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9982'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08825'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.12'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.490'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.003'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.64%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001357'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.741.95'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '96.12'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07208'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.79'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.36'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '25.035.49'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.389'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.968'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.02%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.115'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +13.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '163.03'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '151.26'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +8.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.367'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.623'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +18.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.946.24'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.03179'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +6.38%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.245'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.047'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2903'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.06'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09587'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8318'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.08'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.82%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.17'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.698'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7418'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.261'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.404'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08804'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +7.52%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '140.04'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.57%  '
